$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column (D) values are written as text, matching the
# original inline-string cell type (avoids "214.36" etc. being
# auto-converted to a number by Excel).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.932.72"
$ws.Range("E2").Value = "  +0.81%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.634.26"
$ws.Range("E3").Value = "  +2.05%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.36"
$ws.Range("E5").Value = "  +0.92%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.518"
$ws.Range("E6").Value = "  +0.27%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "28.53"
$ws.Range("E8").Value = "  +1.41%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.258"
$ws.Range("E9").Value = "  +1.47%  "

$ws.Range("E10").Value = "  +0.84%  "

$ws.Range("E11").Value = "  +0.25%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.867.01"
$ws.Range("E12").Value = "  +1.92%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.636.33"
$ws.Range("E13").Value = "  +2.35%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.36"
$ws.Range("E15").Value = "  +18.73%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.950.91"
$ws.Range("E16").Value = "  +0.77%  "

$ws.Range("E17").Value = "  +2.45%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.04"
$ws.Range("E18").Value = "  +0.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.15"
$ws.Range("E19").Value = "  +0.21%  "

$ws.Range("E20").Value = "  +0.64%  "

$ws.Range("E21").Value = "  +0.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.87"
$ws.Range("E22").Value = "  +4.82%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.14"
$ws.Range("E23").Value = "  +2.80%  "

$ws.Range("E24").Value = "  +2.37%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.56"
$ws.Range("E25").Value = "  +1.76%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.51"
$ws.Range("E26").Value = "  +0.54%  "

$ws.Range("E27").Value = "  +0.83%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.62"
$ws.Range("E28").Value = "  +2.72%  "

$ws.Range("E29").Value = "  +0.15%  "

$ws.Range("E30").Value = "  +1.97%  "

$ws.Range("E31").Value = "  +4.38%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.38"
$ws.Range("E32").Value = "  +4.34%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.18"
$ws.Range("E33").Value = "  -0.28%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.425.56"
$ws.Range("E34").Value = "  +0.28%  "

$ws.Range("E35").Value = "  +5.19%  "

$ws.Range("E36").Value = "  -0.28%  "

$ws.Range("E37").Value = "  -3.28%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.29"
$ws.Range("E38").Value = "  -0.05%  "

$ws.Range("B39").Value = "Aave"
$ws.Range("C39").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "76.24"
$ws.Range("E39").Value = "  +13.22%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0169"
$ws.Range("E40").Value = "  +0.44%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.551"
$ws.Range("E41").Value = "  +1.32%  "

$ws.Range("E42").Value = "  +2.92%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.827"
$ws.Range("E43").Value = "  +1.47%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0489"
$ws.Range("E44").Value = "  -0.84%  "

$ws.Range("E45").Value = "  +3.53%  "

$ws.Range("B46").Value = "BitcoinSV"
$ws.Range("C46").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "53.04"
$ws.Range("E46").Value = "  -5.15%  "

$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  +0.24%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.36"
$ws.Range("E48").Value = "  -0.22%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.775.40"
$ws.Range("E49").Value = "  +2.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "90.34"
$ws.Range("E50").Value = "  +4.27%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0112"
$ws.Range("E51").Value = "  +8.59%  "
